# Add the "2022-Q3" data sheet (right after the "总计" summary sheet) and
# update the summary sheet with the new quarter's totals.
#
# Resulting sheet order: 总计, 2022-Q3, 2022-Q2, 2022-Q1, 2021-Q4, 2021-Q3

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert a brand new worksheet right after "总计" and name it "2022-Q3"
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)
$q3 = $wb.Worksheets.Add($null, $summary)
$q3.Name = "2022-Q3"

# Columns B (fund code) and D:G (numeric-looking figures kept as text, same
# as the other quarterly sheets) must be forced to Text format BEFORE the
# values are written, otherwise leading zeros get stripped / the values get
# silently re-typed as numbers.
$q3.Range("B2:B5").NumberFormat = "@"
$q3.Range("D2:G5").NumberFormat = "@"

# Header row (same columns/order used by the other quarterly sheets)
$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

$hdr = $q3.Range("B1:H1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160
$hdr.Borders.LineStyle = 1

# Data rows
$q3.Range("A2").Value = 0
$q3.Range("B2").Value = "011081"
$q3.Range("C2").Value = "国投瑞银港股通混合C"
$q3.Range("D2").Value = "19.35"
$q3.Range("E2").Value = "90.31"
$q3.Range("F2").Value = "2.58"
$q3.Range("G2").Value = "0.4992"
$q3.Range("H2").Value = 10

$q3.Range("A3").Value = 1
$q3.Range("B3").Value = "007110"
$q3.Range("C3").Value = "国投瑞银港股通价值发现混合"
$q3.Range("D3").Value = "19.35"
$q3.Range("E3").Value = "90.31"
$q3.Range("F3").Value = "2.58"
$q3.Range("G3").Value = "0.4992"
$q3.Range("H3").Value = 10

$q3.Range("A4").Value = 2
$q3.Range("B4").Value = "013357"
$q3.Range("C4").Value = "大摩沪港深精选混合C"
$q3.Range("D4").Value = "1.53"
$q3.Range("E4").Value = "92.27"
$q3.Range("F4").Value = "5.20"
$q3.Range("G4").Value = "0.0796"
$q3.Range("H4").Value = 10

$q3.Range("A5").Value = 3
$q3.Range("B5").Value = "013356"
$q3.Range("C5").Value = "大摩沪港深精选混合A"
$q3.Range("D5").Value = "0.68"
$q3.Range("E5").Value = "92.27"
$q3.Range("F5").Value = "5.20"
$q3.Range("G5").Value = "0.0354"
$q3.Range("H5").Value = 10

$idxCol = $q3.Range("A2:A5")
$idxCol.Font.Bold = $true
$idxCol.HorizontalAlignment = -4108
$idxCol.VerticalAlignment = -4160
$idxCol.Borders.LineStyle = 1

# ---------------------------------------------------------------------
# 2. Update the "总计" summary sheet: push every existing row down by one
#    and insert the new "2022-Q3" totals at the top of the data.
# ---------------------------------------------------------------------
$summary.Range("A6").Value = 4
$summary.Range("B6").Value = "2021-Q3"
$summary.Range("C6").Value = 4
$summary.Range("D6").Value = 0.74

$summary.Range("A5").Value = 3
$summary.Range("B5").Value = "2021-Q4"
$summary.Range("C5").Value = 2
$summary.Range("D5").Value = 1.43

$summary.Range("A4").Value = 2
$summary.Range("B4").Value = "2022-Q1"
$summary.Range("C4").Value = 6
$summary.Range("D4").Value = 1.93

$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2022-Q2"
$summary.Range("C3").Value = 2
$summary.Range("D3").Value = 1.36

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 4
$summary.Range("D2").Value = 1.11

# Row 6 is brand new territory on this sheet (previously the table stopped at
# row 5), so its index cell needs the same bold/centered/bordered look as the
# other index cells (A2:A5) in column A.
$newIdxCell = $summary.Range("A6")
$newIdxCell.Font.Bold = $true
$newIdxCell.HorizontalAlignment = -4108
$newIdxCell.VerticalAlignment = -4160
$newIdxCell.Borders.LineStyle = 1

# Keep the originally-active tab ("2021-Q3", now the last sheet) selected,
# since adding a worksheet otherwise shifts the active tab to the new sheet.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$lastSheet.Activate()

Write-Host "2022-Q3 sheet added; summary sheet updated"
